# Add the "(Source3)" data-source suffix to the header columns of both the
# normal layout (row 1) and the interlaced/key layout (row 11). Columns
# B/E, C/F and D/G each reuse the same shared header text, so every one of
# them has to be written explicitly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "FIRST_NAME (Source3)"
$ws.Range("C1").Value = "LAST_NAME (Source3)"
$ws.Range("D1").Value = "AGE (Source3)"
$ws.Range("E1").Value = "FIRST_NAME (Source3)"
$ws.Range("F1").Value = "LAST_NAME (Source3)"
$ws.Range("G1").Value = "AGE (Source3)"

$ws.Range("B11").Value = "FIRST_NAME (Source3)"
$ws.Range("C11").Value = "LAST_NAME (Source3)"
$ws.Range("D11").Value = "AGE (Source3)"
$ws.Range("E11").Value = "FIRST_NAME (Source3)"
$ws.Range("F11").Value = "LAST_NAME (Source3)"
$ws.Range("G11").Value = "AGE (Source3)"

# The header columns use bestFit/autofit column widths, so widening the
# header text (10 extra characters: " (Source3)") widens the columns B-G.
# ColumnWidth values below were picked so the saved <col width="..."> lands
# on (i.e. as close as achievable to) the recorded widths: ~26.105 for the
# FIRST_NAME columns, ~25.598 for LAST_NAME, ~18.391 for AGE.
$ws.Columns.Item(2).ColumnWidth = 25.333333333333332
$ws.Columns.Item(3).ColumnWidth = 24.833333333333332
$ws.Columns.Item(4).ColumnWidth = 17.5
$ws.Columns.Item(5).ColumnWidth = 25.333333333333332
$ws.Columns.Item(6).ColumnWidth = 24.833333333333332
$ws.Columns.Item(7).ColumnWidth = 17.5
